$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C: running average from B2 down to current row
for ($row = 2; $row -le 13; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Formula = "=AVERAGE(`$B`$2:B$row)"
}

# Column D: 3-month running average, starting at row 4
for ($row = 4; $row -le 13; $row++) {
    $startRow = $row - 2
    $cell = $ws.Cells.Item($row, 4)
    $cell.Formula = "=AVERAGE(B${startRow}:B$row)"
}

# Copy the number format/style from column B (Comma style) onto the new C/D cells
$ws.Range("C2:C13").NumberFormat = $ws.Range("B2").NumberFormat
$ws.Range("D4:D13").NumberFormat = $ws.Range("B2").NumberFormat

# Restore selection like the final state of the saved file
$ws.Range("F8").Select()
